# Apply "added submiting to form template" edit:
#  - replace the placeholder header/sample rows with the real scouting
#    form header (matchNum / TeamNum / climb / throw / additionalNotes)
#  - add a new "additionalNotes" column (E)
#  - remove the old sample data rows (3-4) and add real submissions in
#    rows 10-13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text even when it looks numeric
# (plain .Value assignment of a purely-numeric string gets reinterpreted
# as a number, so we go through a text formula and then freeze it back
# down to a literal value with paste-special).
function Set-TextValue {
    param($Range, $Text)
    $escaped = $Text.Replace('"', '""')
    $Range.Formula = '="' + $escaped + '"'
    $Range.Copy()
    $Range.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false

# Clear the old sample rows (3 and 4) entirely
$ws.Range("A3:D4").ClearContents()

# --- Header rows (1 and 2) ---
$headers = @("matchNum", "TeamNum", "climb", "throw", "additionalNotes")
$headerCols = @("A", "B", "C", "D", "E")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headers[$i]
    $ws.Range($headerCols[$i] + "2").Value = $headers[$i]
}

# --- Submitted data rows (10-13) ---
Set-TextValue $ws.Range("A10") "324"
Set-TextValue $ws.Range("B10") "1234"
$ws.Range("C10").Value = "No"
$ws.Range("D10").Value = "Yes"
$ws.Range("E10").Value = "49iulfkhdjkhlksajd23"

Set-TextValue $ws.Range("A11") "341"
Set-TextValue $ws.Range("B11") "2341234"
$ws.Range("C11").Value = "No"
$ws.Range("D11").Value = "Yes"
$ws.Range("E11").Value = "sd32wedaslr2h14lhrkjasgkh"

Set-TextValue $ws.Range("A12") "231"
Set-TextValue $ws.Range("B12") "3214"
$ws.Range("C12").Value = "No"
$ws.Range("D12").Value = "Yes"
$ws.Range("E12").Value = "hdfkjhlskdf"

Set-TextValue $ws.Range("A13") "53"
Set-TextValue $ws.Range("B13") "4123"
$ws.Range("C13").Value = "Yes"
$ws.Range("D13").Value = "No"
Set-TextValue $ws.Range("E13") "421341234"

$excel.CutCopyMode = $false
